$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 509.36365
$ws.Range("J53").Value = 574.625
$ws.Range("L53").Value = 574.625
$ws.Range("N53").Value = -1848.625
$ws.Range("H64").Value = 58330.832
$ws.Range("I64").Value = 127545
$ws.Range("J64").Value = 2959.5
$ws.Range("K64").Value = 127545
$ws.Range("L64").Value = 2959.5
$ws.Range("M64").Value = -127297
$ws.Range("N64").Value = -3455.5
$ws.Range("H67").Value = 58330.832
$ws.Range("I67").Value = 127545
$ws.Range("J67").Value = 2959.5
$ws.Range("K67").Value = 127545
$ws.Range("L67").Value = 2959.5
$ws.Range("M67").Value = -126687
$ws.Range("N67").Value = -4675.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 128.2
$ws.Range("I4").Value = 95.5
$ws.Range("K4").Value = 95.5
$ws.Range("M4").Value = 20.5
$ws.Range("H32").Value = 25671.904
$ws.Range("I32").Value = 7202.053
$ws.Range("J32").Value = 198826.75
$ws.Range("K32").Value = 7202.053
$ws.Range("L32").Value = 198826.75
$ws.Range("M32").Value = -6915.053
$ws.Range("N32").Value = -199400.75
$ws.Range("H33").Value = 14998
$ws.Range("I33").Value = 14998
$ws.Range("K33").Value = 14998
$ws.Range("M33").Value = -14669
$ws.Range("H36").Value = 776
$ws.Range("I36").Value = 776
$ws.Range("K36").Value = 776
$ws.Range("M36").Value = -430
$ws.Range("H102").Value = 73832
$ws.Range("I102").Value = 102346.8
$ws.Range("K102").Value = 102346.8
$ws.Range("M102").Value = -100724.8
$ws.Range("H122").Value = 2032.3572
$ws.Range("I122").Value = 2168.5454
$ws.Range("J122").Value = 1533
$ws.Range("K122").Value = 6505.6362
$ws.Range("L122").Value = 4599
$ws.Range("M122").Value = -4055.6362
$ws.Range("N122").Value = -9499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1691.9714
$ws.Range("I99").Value = 1143.579
$ws.Range("J99").Value = 2343.1875
$ws.Range("K99").Value = 1143.579
$ws.Range("L99").Value = 2343.1875
$ws.Range("M99").Value = 354.421
$ws.Range("N99").Value = -5339.1875
$ws.Range("H105").Value = 135177.33
$ws.Range("I105").Value = 78880.69500000001
$ws.Range("J105").Value = 501105.5
$ws.Range("K105").Value = 78880.69500000001
$ws.Range("L105").Value = 501105.5
$ws.Range("M105").Value = -77133.69500000001
$ws.Range("N105").Value = -504599.5
$ws.Range("H134").Value = 12704.739
$ws.Range("I134").Value = 14176
$ws.Range("J134").Value = 4507.7144
$ws.Range("K134").Value = 42528
$ws.Range("L134").Value = 13523.1432
$ws.Range("M134").Value = -39993
$ws.Range("N134").Value = -18593.1432
$ws.Range("H140").Value = 45080
$ws.Range("J140").Value = 45080
$ws.Range("L140").Value = 45080
$ws.Range("N140").Value = -55440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 3566.6667
$ws.Range("I12").Value = 350
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 350
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -180
$ws.Range("N12").Value = -10340
$ws.Range("H33").Value = 3530.4
$ws.Range("I33").Value = 3530.4
$ws.Range("K33").Value = 3530.4
$ws.Range("M33").Value = -3151.4
$ws.Range("H132").Value = 3072.5
$ws.Range("I132").Value = 3326.2632
$ws.Range("J132").Value = 2536.7778
$ws.Range("K132").Value = 9978.7896
$ws.Range("L132").Value = 7610.3334
$ws.Range("M132").Value = -7448.7896
$ws.Range("N132").Value = -12670.3334
$ws.Range("H134").Value = 1566.7084
$ws.Range("I134").Value = 1224.5
$ws.Range("J134").Value = 2593.3333
$ws.Range("K134").Value = 3673.5
$ws.Range("L134").Value = 7779.999899999999
$ws.Range("M134").Value = -1138.5
$ws.Range("N134").Value = -12849.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1400
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H127").Value = 1131.125
$ws.Range("J127").Value = 1131.125
$ws.Range("L127").Value = 3393.375
$ws.Range("N127").Value = -13313.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 24900
$ws.Range("J58").Value = 24900
$ws.Range("L58").Value = 24900
$ws.Range("N58").Value = -25454
$ws.Range("H102").Value = 224935.67
$ws.Range("I102").Value = 1504.8
$ws.Range("K102").Value = 1504.8
$ws.Range("M102").Value = 117.2
$ws.Range("H126").Value = 5350492
$ws.Range("I126").Value = 3168.8572
$ws.Range("J126").Value = 14708307
$ws.Range("K126").Value = 9506.571599999999
$ws.Range("L126").Value = 44124921
$ws.Range("M126").Value = -7036.571599999999
$ws.Range("N126").Value = -44129861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4300
$ws.Range("H22").Value = 1776.875
$ws.Range("J22").Value = 1448.5
$ws.Range("L22").Value = 1448.5
$ws.Range("N22").Value = -2038.5
$ws.Range("H27").Value = 1776.875
$ws.Range("J27").Value = 1448.5
$ws.Range("L27").Value = 1448.5
$ws.Range("N27").Value = -1662.5
$ws.Range("H40").Value = 61216.707
$ws.Range("I40").Value = 144615.28
$ws.Range("J40").Value = 2837.7
$ws.Range("K40").Value = 144615.28
$ws.Range("L40").Value = 2837.7
$ws.Range("M40").Value = -144479.28
$ws.Range("N40").Value = -3109.7
$ws.Range("H61").Value = 1334.6072
$ws.Range("I61").Value = 1280.9474
$ws.Range("J61").Value = 1447.8889
$ws.Range("K61").Value = 1280.9474
$ws.Range("L61").Value = 1447.8889
$ws.Range("M61").Value = -1078.9474
$ws.Range("N61").Value = -1851.8889
$ws.Range("H68").Value = 4294
$ws.Range("J68").Value = 4618.25
$ws.Range("L68").Value = 4618.25
$ws.Range("N68").Value = -6116.25
$ws.Range("H71").Value = 4294
$ws.Range("J71").Value = 4618.25
$ws.Range("L71").Value = 23091.25
$ws.Range("N71").Value = -30579.25
$ws.Range("H82").Value = 1793.4667
$ws.Range("J82").Value = 2804
$ws.Range("L82").Value = 2804
$ws.Range("N82").Value = -3526
$ws.Range("H85").Value = 1793.4667
$ws.Range("J85").Value = 2804
$ws.Range("L85").Value = 2804
$ws.Range("N85").Value = -5300
$ws.Range("H93").Value = 2311.6538
$ws.Range("I93").Value = 2315.6667
$ws.Range("J93").Value = 2306.182
$ws.Range("K93").Value = 2315.6667
$ws.Range("L93").Value = 2306.182
$ws.Range("M93").Value = -1067.6667
$ws.Range("N93").Value = -4802.182
$ws.Range("H100").Value = 1964.1666
$ws.Range("I100").Value = 1750
$ws.Range("J100").Value = 2392.5
$ws.Range("K100").Value = 1750
$ws.Range("L100").Value = 2392.5
$ws.Range("M100").Value = -1209
$ws.Range("N100").Value = -3474.5
$ws.Range("H113").Value = 1334.6072
$ws.Range("I113").Value = 1280.9474
$ws.Range("J113").Value = 1447.8889
$ws.Range("K113").Value = 1280.9474
$ws.Range("L113").Value = 1447.8889
$ws.Range("M113").Value = 889.0526
$ws.Range("N113").Value = -5787.8889
$ws.Range("H122").Value = 2641.3572
$ws.Range("I122").Value = 2641.3572
$ws.Range("K122").Value = 7924.071599999999
$ws.Range("M122").Value = -5474.071599999999
$ws.Range("H126").Value = 4300
$ws.Range("H141").Value = 52753.75
$ws.Range("J141").Value = 52753.75
$ws.Range("L141").Value = 52753.75
$ws.Range("N141").Value = -63113.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 18200
$ws.Range("J40").Value = 18200
$ws.Range("L40").Value = 18200
$ws.Range("N40").Value = -18498
$ws.Range("H56").Value = 26106
$ws.Range("I56").Value = 3800
$ws.Range("J56").Value = 40976.668
$ws.Range("K56").Value = 3800
$ws.Range("L56").Value = 40976.668
$ws.Range("M56").Value = -3086
$ws.Range("N56").Value = -42404.668
$ws.Range("H126").Value = 1291.0454
$ws.Range("I126").Value = 1426.25
$ws.Range("J126").Value = 930.5
$ws.Range("K126").Value = 4278.75
$ws.Range("L126").Value = 2791.5
$ws.Range("M126").Value = -1808.75
$ws.Range("N126").Value = -7731.5
$ws.Range("H132").Value = 3267.4856
$ws.Range("I132").Value = 3452
$ws.Range("J132").Value = 2644.75
$ws.Range("K132").Value = 10356
$ws.Range("L132").Value = 7934.25
$ws.Range("M132").Value = -7826
$ws.Range("N132").Value = -12994.25
